$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update "Compilation success" row (row 5)
$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = ""

# Update "Runtime without error" row (row 6)
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "Missing initial redirect"

# Update Code BLEU score (row 12)
$ws.Range("B12").Value = 0.2924470879319078
$ws.Range("C12").Value = "{'codebleu': 0.29244708793190777, 'ngram_match_score': 0.14075717909287128, 'weighted_ngram_match_score': 0.15906780267138987, 'syntax_match_score': 0.5604395604395604, 'dataflow_match_score': 0.30952380952380953}"

# Update the active selection to B7, mirroring the saved selection in the workbook
$ws.Range("B7").Select()
